$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph that currently follows
# the H1 title ("Play Country Farming Slot for Free - Review 2021"). That
# whole paragraph (bold label "Meta description" + the description text) is
# deleted outright.
$d.Paragraphs.Item(2).Range.Delete() | Out-Null

# --- Step 2: near the end of the document, insert a brand-new bold
# paragraph reading "Play Country Farming Slot for Free - Review 2021"
# right before the final ("Prompt: ...") paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$newParaIndex = $d.Paragraphs.Count - 1
$newParaRange = $d.Paragraphs.Item($newParaIndex).Range
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Country Farming Slot for Free - Review 2021</w:t></w:r></w:p>'
$newParaRange.InsertXML($newParaXml) | Out-Null

# --- Step 3: the final paragraph (previously the "Prompt: ..." image-prompt
# text, still italic) now becomes the old meta-description copy instead.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Find.Execute(
    "Prompt: Create a feature image for Country Farming that showcases the game's fun and colorful nature while incorporating the Maya warrior character with glasses. The image should be in a cartoon style and feature the game's farm symbols such as animals and fruits. It should also include the game's logo. The Maya warrior can be depicted engaging in a fun activity or interacting with the symbols in some way, making the image lively and engaging to potential players.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Our review of Country Farming slot - a fun and engaging video game with high winning potential. Play now for free!",
    2) | Out-Null
